$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values on existing rows 3 and 4
$ws.Range("A3").Value = 12
$ws.Range("A4").Value = 9

# Add new row 5, replicating the formatting of row 4 (same style, same grade text)
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "1E"

# Update the active selection to A2
$ws.Range("A2").Select()
